$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 5) to the sheet, mirroring the layout/formatting of row 4.
$ws.Range("A4:I4").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Value = 42647.680567129632
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = 9988.1
$ws.Range("D5").Value = 10033.25
$ws.Range("E5").Value = 313
$ws.Range("F5").Value = 311.58999999999997
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = -0.45
$ws.Range("I5").Value = $true
